$d = $word.ActiveDocument

function Replace-Exact($oldText, $newText) {
    $text = $d.Content.Text
    $idx = $text.IndexOf($oldText)
    if ($idx -lt 0) {
        throw "Could not find text: $oldText"
    }
    $rng = $d.Range($idx, $idx + $oldText.Length)
    if ($rng.Text -ne $oldText) {
        throw "Range text mismatch for: $oldText"
    }
    $rng.Text = $newText
}

# 1. "data " -> "sandwiches " (attach line context)
Replace-Exact "data " "sandwiches "

# 2. "(data) " -> "(sandwiches) "
Replace-Exact "(data) " "(sandwiches) "

# 3. Bootstrap sample listing: rewrite the 12-line verbatim block into an 8-line block.
Replace-Exact '##  [1] "Prawn, Mayo"                  "Tuna, Cucumber"              ' '##  [1] "Egg, Bacon"           "Egg, Bacon"           "Sausage, Brown Sauce"'
Replace-Exact '##  [3] "Chicken, Sweetcorn"           "Chicken, Bacon"              ' '##  [4] "Ham, Mayo"            "Ham, Mustard"         "Chicken Salad"       '
Replace-Exact '##  [5] "Ham, Salad"                   "Ham, Mustard"                ' '##  [7] "Cheese, Pickle"       "Cheese Ploughman "    "Ham, Mustard"        '
Replace-Exact '##  [7] "Sausage, Brown Sauce"         "Tuna, Cucumber"              ' '## [10] "Ham, Mustard"         "Ham, Mayo"            "Sausage, Brown Sauce"'
Replace-Exact '##  [9] "Prawn, Mayo"                  "Ham, Cheese"                 ' '## [13] "Chicken, Stuffing"    "Chicken Salad"        "Tuna, Cucumber"      '
Replace-Exact '## [11] "Prawn, Mayo"                  "Cheese, Tomato"              ' '## [16] "Cheese, Mayo"         "Ham, Mustard"         "Ham, Egg"            '
Replace-Exact '## [13] "Tuna, Cucumber"               "Ham, Egg"                    ' '## [19] "Egg, Rocket"          "Chicken Salad"        "Ham, Mustard"        '
Replace-Exact '## [15] "Prawn, Mayo"                  "Chicken, Bacon"              ' '## [22] "Breakfast"            "Chicken, Bacon"       "Breakfast"'

# Now delete the now-redundant trailing lines 17/19/21/23 (with their line breaks).
$text = $d.Content.Text
$afterNewLine15 = $text.IndexOf('## [22] "Breakfast"            "Chicken, Bacon"       "Breakfast"') + '## [22] "Breakfast"            "Chicken, Bacon"       "Breakfast"'.Length
$line23 = '## [23] "Ham, Mayo"                    "Ham, Egg"'
$idxEnd = $text.IndexOf($line23) + $line23.Length
$rng = $d.Range($afterNewLine15, $idxEnd)
$rng.Delete()

Write-Output "done-part1"
